$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Platform Coverage")

# Update row 2 coverage values (H2:AD2) from 0.6 to 0.736
$ws.Range("H2:AD2").Value = 0.736

# Add new row 12: "Vector Control" label in B12, with two tiny coverage
# values (1E-8) in J12 and K12, styled with a black-colored font.
$ws.Range("B12").Value = "Vector Control"
$ws.Range("J12:K12").Value = 0.00000001
$ws.Range("J12:K12").Font.Color = 0

# Match the saved selection state on the sheet
[void]$ws.Range("B12:K12").Select()

Write-Host "edit applied"
